$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before column N ---
# This shifts the old N/O/P columns (Late / Outstanding / Disbursement) one
# column to the right, becoming O/P/Q, and leaves a new blank column N in
# between "In Advance" (M) and "Late" (now O).
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N:N").Insert()
# Match the width of the new blank column to its left neighbour (column M),
# same as Excel typically carries over formatting on a manual insert.
$wsRepay.Columns("N:N").ColumnWidth = $wsRepay.Columns("M:M").ColumnWidth
[void]$wsRepay.Range("S7").Select()

# --- Summary sheet becomes the active/selected sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
[void]$wsSummary.Range("H9").Select()
